$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Cluster 3)
$ws.Range("E4").Value = 0.1983471074380164

# Row 5 (Cluster 4)
$ws.Range("C5").Value = 0.368421052631579
$ws.Range("D5").Value = 0.02631578947368422
$ws.Range("E5").Value = 0.02631578947368423
$ws.Range("G5").Value = 0.02631578947368423
$ws.Range("I5").Value = 0.1842105263157895
$ws.Range("J5").Value = 0.2105263157894737

# Row 6 (Cluster 5)
$ws.Range("B6").Value = 0.03749999999999996
$ws.Range("C6").Value = 0.03749999999999999
$ws.Range("D6").Value = 0.1000000000000001
$ws.Range("E6").Value = 0.1125
$ws.Range("F6").Value = 0.1875
$ws.Range("I6").Value = 0.1625
$ws.Range("J6").Value = 0.04999999999999998

# Row 9 (Cluster 8)
$ws.Range("E9").Value = 0.4358974358974359
$ws.Range("J9").Value = 0.1025641025641025
$ws.Range("K9").Value = 0.05128205128205127

# Row 10 (Cluster 9)
$ws.Range("I10").Value = 0.1090909090909091
$ws.Range("J10").Value = 0.2363636363636364

# Row 11 (Cluster 10)
$ws.Range("C11").Value = 0.3714285714285714
$ws.Range("J11").Value = 0.02857142857142857
